# Re-generate the statistics with the fixed minutes and seconds formatting
# in the "haul" (Общее время) field: zero-pad single-digit minutes and
# seconds, e.g. "185 ч. 45 мин. 8 сек." -> "185 ч. 45 мин. 08 сек."
# Hours are left untouched even when they are a single digit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$hoursCol = 9  # column I = "Общее время"

$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $hoursCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours   = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]

        $needsFix = $false
        if ($minutes.Length -eq 1) { $minutes = "0$minutes"; $needsFix = $true }
        if ($seconds.Length -eq 1) { $seconds = "0$seconds"; $needsFix = $true }

        if ($needsFix) {
            $cell.Value2 = "$hours ч. $minutes мин. $seconds сек."
            $changed++
        }
    }
}

Write-Output "Fixed $changed haul time cell(s)."
